$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.354.27'
$ws.Range('E2').Value = '  +0.80%  '
$ws.Range('D3').Value = '1.687.86'
$ws.Range('E4').Value = '  +0.74%  '
$ws.Range('D5').Value = "'218.39"
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').Value = "'0.5473"
$ws.Range('E6').Value = '  +4.94%  '
$ws.Range('D7').Value = "'1.010"
$ws.Range('E7').Value = '  +0.68%  '
$ws.Range('D8').Value = "'0.2729"
$ws.Range('E8').Value = '  +1.14%  '
$ws.Range('D9').Value = "'0.06456"
$ws.Range('E9').Value = '  +0.92%  '
$ws.Range('E10').Value = '  +1.00%  '
$ws.Range('D11').Value = "'0.07681"
$ws.Range('E11').Value = '  +3.39%  '
$ws.Range('D12').Value = '1.687.98'
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('D13').Value = "'4.536"
$ws.Range('D14').Value = "'0.5813"
$ws.Range('E14').Value = '  -0.28%  '
$ws.Range('D15').Value = "'0.000008325"
$ws.Range('E15').Value = '  -2.30%  '
$ws.Range('D16').Value = "'65.10"
$ws.Range('E16').Value = '  +1.34%  '
$ws.Range('D17').Value = '26.414.07'
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').Value = "'4.942"
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D20').Value = "'10.98"
$ws.Range('E20').Value = '  +1.71%  '
$ws.Range('D21').Value = "'190.41"
$ws.Range('E21').Value = '  +0.39%  '
$ws.Range('D22').Value = "'6.231"
$ws.Range('E23').Value = '  +0.70%  '
$ws.Range('D24').Value = "'149.39"
$ws.Range('E24').Value = '  +3.15%  '
$ws.Range('D25').Value = "'0.1311"
$ws.Range('E25').Value = '  +5.75%  '
$ws.Range('D26').Value = "'7.881"
$ws.Range('E26').Value = '  +3.43%  '
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('D28').Value = "'0.06347"
$ws.Range('E28').Value = '  -3.52%  '
$ws.Range('D29').Value = "'1.413"
$ws.Range('E29').Value = '  +6.36%  '
$ws.Range('D30').Value = "'1.329"
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('D31').Value = "'3.575"
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('D32').Value = "'3.577"
$ws.Range('E33').Value = '  +0.31%  '
$ws.Range('E34').Value = '  +2.44%  '
$ws.Range('D35').Value = "'0.6176"
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('D36').Value = "'2.412"
$ws.Range('E36').Value = '  +1.93%  '
$ws.Range('D37').Value = "'2.719"
$ws.Range('E37').Value = '  +0.75%  '
$ws.Range('D38').Value = "'6.247"
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('D39').Value = '1.115.55'
$ws.Range('E39').Value = '  +1.76%  '
$ws.Range('D40').Value = "'0.01629"
$ws.Range('E40').Value = '  +1.95%  '
$ws.Range('D41').Value = "'0.8766"
$ws.Range('E41').Value = '  +0.51%  '
$ws.Range('D42').Value = "'1.015"
$ws.Range('E42').Value = '  +0.41%  '
$ws.Range('D43').Value = "'100.97"
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('D44').Value = '1.840.10'
$ws.Range('E44').Value = '  +1.07%  '
$ws.Range('D45').Value = "'0.00000000109"
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('D46').Value = "'57.26"
$ws.Range('E46').Value = '  +1.33%  '
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').Value = "'1.014"
$ws.Range('E47').Value = '  +1.14%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = "'8.207"
$ws.Range('E48').Value = '  +0.70%  '
$ws.Range('D49').Value = "'0.05272"
$ws.Range('E49').Value = '  +0.59%  '
$ws.Range('D50').Value = "'0.4303"
$ws.Range('E50').Value = '  +0.57%  '
$ws.Range('E51').Value = '  +0.66%  '
